$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, pushing the former row 100 (and everything
# below it) down by one row. This expands the used range from A1:R123 to
# A1:R124, matching the rest of the data's existing layout/formatting
# (the date column D keeps its "2" style / date number format automatically).
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new record.
$ws.Cells.Item(100, 1).Value = 10
$ws.Cells.Item(100, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(100, 3).Value = "La Araucanía"
$ws.Cells.Item(100, 4).Value = 45015
$ws.Cells.Item(100, 5).Value = 9
$ws.Cells.Item(100, 6).Value = 100112022
$ws.Cells.Item(100, 7).Value = "Arveja Verde"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 40
$ws.Cells.Item(100, 11).Value = 23000
$ws.Cells.Item(100, 12).Value = 23000
$ws.Cells.Item(100, 13).Value = 23000
$ws.Cells.Item(100, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(100, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(100, 16).Value = 920
$ws.Cells.Item(100, 17).Value = 25
$ws.Cells.Item(100, 18).Value = "Hortaliza"
